$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "'56717"
$ws.Range("A2").ClearFormats()
$ws.Range("E2").Value = "'2807"
$ws.Range("E2").ClearFormats()
$ws.Range("A3").Value = "'67692"
$ws.Range("A3").ClearFormats()
$ws.Range("A5").Value = "'50410"
$ws.Range("A5").ClearFormats()
$ws.Range("A7").Value = "'45169"
$ws.Range("A7").ClearFormats()
$ws.Range("A8").Value = "'10683"
$ws.Range("A8").ClearFormats()
$ws.Range("E8").Value = "'5707"
$ws.Range("E8").ClearFormats()
$ws.Range("A9").Value = "'13943"
$ws.Range("A9").ClearFormats()
$ws.Range("E9").Value = "'5440"
$ws.Range("E9").ClearFormats()
$ws.Range("A10").Value = "'17494"
$ws.Range("A10").ClearFormats()
$ws.Range("A11").Value = "'18715"
$ws.Range("A11").ClearFormats()
$ws.Range("A12").Value = "'29166"
$ws.Range("A12").ClearFormats()
$ws.Range("E12").Value = "'4598"
$ws.Range("E12").ClearFormats()
$ws.Range("A13").Value = "'53432"
$ws.Range("A13").ClearFormats()
$ws.Range("A14").Value = "'67524"
$ws.Range("A14").ClearFormats()
$ws.Range("A16").Value = "'14035"
$ws.Range("A16").ClearFormats()
$ws.Range("A17").Value = "'15987"
$ws.Range("A17").ClearFormats()
$ws.Range("B17").Value = "'31495601"
$ws.Range("B17").ClearFormats()
$ws.Range("C17").Value = "陈晓军"
$ws.Range("E17").Value = "'5295"
$ws.Range("E17").ClearFormats()
$ws.Range("A18").Value = "'16083"
$ws.Range("A18").ClearFormats()
$ws.Range("B18").Value = "'31134300"
$ws.Range("B18").ClearFormats()
$ws.Range("C18").Value = "McMaX"
$ws.Range("E18").Value = "'5290"
$ws.Range("E18").ClearFormats()
$ws.Range("A19").Value = "'16743"
$ws.Range("A19").ClearFormats()
$ws.Range("B19").Value = "'54698813"
$ws.Range("B19").ClearFormats()
$ws.Range("C19").Value = "閃亮唐老鴨"
$ws.Range("E19").Value = "'5249"
$ws.Range("E19").ClearFormats()
$ws.Range("A20").Value = "'20605"
$ws.Range("A20").ClearFormats()
$ws.Range("E20").Value = "'5028"
$ws.Range("E20").ClearFormats()
$ws.Range("A21").Value = "'21038"
$ws.Range("A21").ClearFormats()
$ws.Range("E21").Value = "'5004"
$ws.Range("E21").ClearFormats()
$ws.Range("A22").Value = "'23341"
$ws.Range("A22").ClearFormats()
$ws.Range("E22").Value = "'4879"
$ws.Range("E22").ClearFormats()
$ws.Range("A23").Value = "'30299"
$ws.Range("A23").ClearFormats()
$ws.Range("A24").Value = "'31600"
$ws.Range("A24").ClearFormats()
$ws.Range("B24").Value = "'58839983"
$ws.Range("B24").ClearFormats()
$ws.Range("C24").Value = "每逢佳节胖六斤"
$ws.Range("E24").Value = "'4497"
$ws.Range("E24").ClearFormats()
$ws.Range("A25").Value = "'32108"
$ws.Range("A25").ClearFormats()
$ws.Range("B25").Value = "'56732705"
$ws.Range("B25").ClearFormats()
$ws.Range("C25").Value = "时间温柔皆遗憾"
$ws.Range("E25").Value = "'4475"
$ws.Range("E25").ClearFormats()
$ws.Range("A26").Value = "'33535"
$ws.Range("A26").ClearFormats()
$ws.Range("E26").Value = "'4408"
$ws.Range("E26").ClearFormats()
$ws.Range("A27").Value = "'39307"
$ws.Range("A27").ClearFormats()
$ws.Range("E27").Value = "'4124"
$ws.Range("E27").ClearFormats()
$ws.Range("A28").Value = "'39944"
$ws.Range("A28").ClearFormats()
$ws.Range("E28").Value = "'4094"
$ws.Range("E28").ClearFormats()
$ws.Range("A29").Value = "'44294"
$ws.Range("A29").ClearFormats()
$ws.Range("A30").Value = "'6022"
$ws.Range("A30").ClearFormats()
$ws.Range("E30").Value = "'6132"
$ws.Range("E30").ClearFormats()
$ws.Range("A31").Value = "'8166"
$ws.Range("A31").ClearFormats()
$ws.Range("A32").Value = "'11171"
$ws.Range("A32").ClearFormats()
$ws.Range("E32").Value = "'5669"
$ws.Range("E32").ClearFormats()
$ws.Range("A33").Value = "'12152"
$ws.Range("A33").ClearFormats()
$ws.Range("A34").Value = "'12412"
$ws.Range("A34").ClearFormats()
$ws.Range("A35").Value = "'14877"
$ws.Range("A35").ClearFormats()
$ws.Range("E35").Value = "'5367"
$ws.Range("E35").ClearFormats()
$ws.Range("A36").Value = "'18029"
$ws.Range("A36").ClearFormats()
$ws.Range("A37").Value = "'20287"
$ws.Range("A37").ClearFormats()
$ws.Range("E37").Value = "'5046"
$ws.Range("E37").ClearFormats()
$ws.Range("A38").Value = "'28680"
$ws.Range("A38").ClearFormats()
$ws.Range("E38").Value = "'4619"
$ws.Range("E38").ClearFormats()
$ws.Range("A39").Value = "'31219"
$ws.Range("A39").ClearFormats()
$ws.Range("B39").Value = "'47459684"
$ws.Range("B39").ClearFormats()
$ws.Range("C39").Value = "㊥阿闹切克闹"
$ws.Range("E39").Value = "'4513"
$ws.Range("E39").ClearFormats()
$ws.Range("A40").Value = "'31702"
$ws.Range("A40").ClearFormats()
$ws.Range("B40").Value = "'56573048"
$ws.Range("B40").ClearFormats()
$ws.Range("C40").Value = "Xiaotian"
$ws.Range("E40").Value = "'4493"
$ws.Range("E40").ClearFormats()
$ws.Range("A41").Value = "'33042"
$ws.Range("A41").ClearFormats()
$ws.Range("E41").Value = "'4431"
$ws.Range("E41").ClearFormats()
$ws.Range("A42").Value = "'33718"
$ws.Range("A42").ClearFormats()
$ws.Range("E42").Value = "'4399"
$ws.Range("E42").ClearFormats()
$ws.Range("A43").Value = "'36727"
$ws.Range("A43").ClearFormats()
$ws.Range("E43").Value = "'4254"
$ws.Range("E43").ClearFormats()
$ws.Range("A44").Value = "'39626"
$ws.Range("A44").ClearFormats()
$ws.Range("A45").Value = "'40682"
$ws.Range("A45").ClearFormats()
$ws.Range("A46").Value = "'42068"
$ws.Range("A46").ClearFormats()
$ws.Range("A47").Value = "'42558"
$ws.Range("A47").ClearFormats()
$ws.Range("E47").Value = "'3953"
$ws.Range("E47").ClearFormats()
$ws.Range("A48").Value = "'42899"
$ws.Range("A48").ClearFormats()
$ws.Range("A49").Value = "'49208"
$ws.Range("A49").ClearFormats()
$ws.Range("A50").Value = "'57184"
$ws.Range("A50").ClearFormats()
$ws.Range("E50").Value = "'2790"
$ws.Range("E50").ClearFormats()
$ws.Range("A51").Value = "'67400"
$ws.Range("A51").ClearFormats()
$ws.Range("A52").Value = "'61723"
$ws.Range("A52").ClearFormats()
$ws.Range("E52").Value = "'2648"
$ws.Range("E52").ClearFormats()
$ws.Range("A53").Value = "'50581"
$ws.Range("A53").ClearFormats()
$ws.Range("A56").Value = "'42499"
$ws.Range("A56").ClearFormats()
$ws.Range("A57").Value = "'51108"
$ws.Range("A57").ClearFormats()
$ws.Range("A58").Value = "'58860"
$ws.Range("A58").ClearFormats()
$ws.Range("E58").Value = "'2732"
$ws.Range("E58").ClearFormats()
$ws.Range("A59").Value = "'68060"
$ws.Range("A59").ClearFormats()
$ws.Range("A60").Value = "'105207"
$ws.Range("A60").ClearFormats()
$ws.Range("A61").Value = "'106852"
$ws.Range("A61").ClearFormats()
$ws.Range("A62").Value = "'109737"
$ws.Range("A62").ClearFormats()
$ws.Range("A63").Value = "'121250"
$ws.Range("A63").ClearFormats()
$ws.Range("A75").Value = "'46939"
$ws.Range("A75").ClearFormats()
$ws.Range("E75").Value = "'3431"
$ws.Range("E75").ClearFormats()
$ws.Range("A78").Value = "'89649"
$ws.Range("A78").ClearFormats()
$ws.Range("E78").Value = "'1904"
$ws.Range("E78").ClearFormats()
$ws.Range("A79").Value = "'96375"
$ws.Range("A79").ClearFormats()
$ws.Range("A80").Value = "'158555"
$ws.Range("A80").ClearFormats()
$ws.Range("A81").Value = "'210794"
$ws.Range("A81").ClearFormats()
